$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 is the "2022 global" (Bolton & Kacperczyk) literature-review entry ---

# Findings (column N) previously had no note for this row; add one.
$ws.Range("N3").Value2 = "Carbon emissions levels are a `npersistent characteristic when `nthe news effect is taken `nout`n- backwards imputing emissions leads to similar premia"
# Match the wrap/top-aligned text formatting used by the rest of the row.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Applicable to my content (column Q): append two more applicability notes.
$existingQ3 = $ws.Range("Q3").Value2
$ws.Range("Q3").Value2 = $existingQ3 + "`n`nThere are horizon effects to consider when modelling`n`nInformation effects analyses of Trucost are likely (almost directly) applicable in my case"

# Scroll the sheet view up so row 2 (instead of row 3) is the top visible row.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 13
